# Add two new rows to the "Sheet2" worksheet (the active sheet) documenting
# how to set the vimrc tab-to-spaces config and the default colorscheme.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A/B for both rows first, then column C bottom row before top
# row -- mirrors the author's original data-entry order so shared-string
# allocation matches.
$ws.Range("A9").Value = "config"
$ws.Range("B9").Value = "set tab to 4 spaces"

$ws.Range("A10").Value = "config"
$ws.Range("B10").Value = "set default colorscheme"

$ws.Range("C10").Value = "In .vimrc, add:`ncolorscheme darkblue      `" myles put 2017-04-07 on set default colorscheme"
$ws.Range("C9").Value = "In .vimrc, add:`nfiletype plugin indent on       `nset tabstop=4                   `"show existing tab with 4 spaces width`nset shiftwidth=4                `" when indenting with '>', use 4 spaces width`nset expandtab                   `" On pressing tab, insert 4 spaces"

# Column C entries wrap text, matching the style used for the other
# multi-line description cells in this sheet.
$ws.Range("C9").WrapText = $true
$ws.Range("C10").WrapText = $true

# Row heights to fit the wrapped content.
$ws.Rows.Item(9).RowHeight = 75
$ws.Rows.Item(10).RowHeight = 30

# Restore the selection to where the author last left it.
$ws.Range("C7").Select() | Out-Null
